$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

# New Cypher query text placed in A2 (new shared string entry); the existing
# cell already carries the wrap-text style (s="1").
$query = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Prostate cancer, NOS''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'
$ws.Range("A2").Value = $query

# Row 2 grows tall enough to show the wrapped query text.
$ws.Rows.Item(2).RowHeight = 87

# Selection moves off the header area and onto the query column.
$ws.Range("B2:B5").Select() | Out-Null

$wb.Save() | Out-Null
